# Updates cryptos list: refresh Price (D) and Volume(1h) (E) columns for rows 2-51
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.471.98"
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("D3").Value = "1.955.57"
$ws.Range("E3").Value = "  +0.64%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.37"
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.628"
$ws.Range("E6").Value = "  +2.58%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.26"
$ws.Range("E7").Value = "  +4.94%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.377"
$ws.Range("E9").Value = "  +3.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0787"
$ws.Range("E10").Value = "  -2.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.104"
$ws.Range("E11").Value = "  +0.82%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.15"
$ws.Range("E12").Value = "  +6.04%  "
$ws.Range("E13").Value = "  +4.07%  "
$ws.Range("D14").Value = "2.238.84"
$ws.Range("E14").Value = "  +0.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.51"
$ws.Range("E15").Value = "  -1.44%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.27"
$ws.Range("E16").Value = "  +1.60%  "
$ws.Range("D17").Value = "1.956.99"
$ws.Range("E17").Value = "  +0.61%  "
$ws.Range("D18").Value = "36.420.37"
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.10"
$ws.Range("E19").Value = "  -0.24%  "
$ws.Range("D20").Value = "0.0₃0853"
$ws.Range("E20").Value = "  -0.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "228.83"
$ws.Range("E21").Value = "  +0.50%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.07"
$ws.Range("E22").Value = "  +2.24%  "
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.44"
$ws.Range("E24").Value = "  +2.01%  "
$ws.Range("E25").Value = "  +2.82%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.143"
$ws.Range("E26").Value = "  +6.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.15"
$ws.Range("E27").Value = "  -0.33%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "160.49"
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.22"
$ws.Range("E29").Value = "  +0.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.30"
$ws.Range("E31").Value = "  +1.98%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.78"
$ws.Range("E32").Value = "  +2.91%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0611"
$ws.Range("E33").Value = "  -0.67%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.45"
$ws.Range("E34").Value = "  +6.78%  "
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.26"
$ws.Range("E36").Value = "  +2.85%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.38"
$ws.Range("E37").Value = "  +4.01%  "
$ws.Range("E38").Value = "  -0.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.42"
$ws.Range("E39").Value = "  -11.82%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0964"
$ws.Range("E40").Value = "  -2.44%  "
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("E42").Value = "  +1.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0210"
$ws.Range("E43").Value = "  +0.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.84"
$ws.Range("E44").Value = "  +0.52%  "
$ws.Range("D45").Value = "1.361.04"
$ws.Range("E45").Value = "  +1.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "88.65"
$ws.Range("E46").Value = "  +2.80%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.03"
$ws.Range("E47").Value = "  +0.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.17"
$ws.Range("E48").Value = "  +0.78%  "
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "46.23"
$ws.Range("E50").Value = "  +6.98%  "
$ws.Range("D51").Value = "2.134.41"
$ws.Range("E51").Value = "  +0.76%  "
